$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1425
$ws.Range("F3").Value = 104
$ws.Range("F5").Value = 6737
$ws.Range("F6").Value = 529
$ws.Range("F8").Value = 39
$ws.Range("F9").Value = 4576
$ws.Range("F10").Value = 6793
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = 227
$ws.Range("F13").Value = 1394
$ws.Range("F14").Value = 804
$ws.Range("F15").Value = 116
$ws.Range("F17").Value = 33
$ws.Range("F18").Value = 1130
$ws.Range("F20").Value = 129
$ws.Range("F22").Value = 186
$ws.Range("F23").Value = 26
$ws.Range("F24").Value = 1058
$ws.Range("F28").Value = 119
$ws.Range("F32").Value = 98
$ws.Range("F37").Value = 515
$ws.Range("F38").Value = 366
$ws.Range("F39").Value = 41
$ws.Range("F40").Value = 52
$ws.Range("F41").Value = 314
$ws.Range("F43").Value = 523
$ws.Range("F44").Value = 63

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 15
$ws.Range("F3").Value = 15
$ws.Range("F11").Value = 31
$ws.Range("F13").Value = 17
$ws.Range("F17").Value = 1727
$ws.Range("F27").Value = 613
$ws.Range("F28").Value = 47
$ws.Range("F31").Value = 772
$ws.Range("F32").Value = 964
$ws.Range("F39").Value = 129
$ws.Range("F42").Value = 61

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 843
$ws.Range("F6").Value = 611
$ws.Range("F8").Value = 1284
$ws.Range("F9").Value = 1832

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 15
$ws.Range("F4").Value = 1425
$ws.Range("F6").Value = 843
$ws.Range("F8").Value = 104
$ws.Range("F9").Value = 611
$ws.Range("F10").Value = 611
$ws.Range("F12").Value = 6737
$ws.Range("F13").Value = 529
$ws.Range("F15").Value = 39
$ws.Range("F16").Value = 4576
$ws.Range("F17").Value = 31
$ws.Range("F18").Value = 6793
$ws.Range("F19").Value = 227
$ws.Range("F20").Value = 1394
$ws.Range("F22").Value = 804
$ws.Range("F23").Value = 116
$ws.Range("F24").Value = 1284
$ws.Range("F26").Value = 1130
$ws.Range("F27").Value = 129
$ws.Range("F28").Value = 186
$ws.Range("F29").Value = 1058
$ws.Range("F30").Value = 613
$ws.Range("F34").Value = 119
$ws.Range("F36").Value = 98
$ws.Range("F39").Value = 964
$ws.Range("F40").Value = 515
$ws.Range("F42").Value = 366
$ws.Range("F43").Value = 41
$ws.Range("F45").Value = 314
$ws.Range("F46").Value = 523
$ws.Range("F50").Value = 61
